# Trade #27 closed at 2026-02-17 08:03:28 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook:
#   - Summary sheet: refresh current capital / P&L / trade-count metrics
#   - Strategy Status sheet: refresh the MarketMaking strategy row
#   - All Trades / MarketMaking sheets: append the newly closed trade (#27)

$wb = $excel.ActiveWorkbook

# --- Summary ---------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.42   # Current Capital
$summary.Range("B4").Value = -0.58     # Total P&L $
$summary.Range("B5").Value = -0.43     # Total P&L %
$summary.Range("B6").Value = 27        # Total Trades
$summary.Range("B8").Value = 13        # Losing Trades
$summary.Range("B9").Value = 25.93     # Win Rate %

# --- Strategy Status (MarketMaking row) -------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.42      # Capital
$status.Range("D4").Value = 27         # Trades
$status.Range("E4").Value = -0.58      # P&L $
$status.Range("F4").Value = -0.58      # P&L %
$status.Range("G4").Value = 25.93      # Win Rate %

# --- Append new closed trade (row 28) to the trade logs ---------------
$tradeSheets = @("All Trades", "MarketMaking")

foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 28

    $ws.Cells.Item($row, 1).Value = 27
    # Leading apostrophe forces the date-like text to stay literal text
    # instead of being auto-converted into a date serial number.
    $ws.Cells.Item($row, 2).Value = "'2026-02-17"
    $ws.Cells.Item($row, 3).Value = "08:03:22"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.4
    $ws.Cells.Item($row, 7).Value = 0.278269
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -30.4327
    $ws.Cells.Item($row, 10).Value = -0.12
    $ws.Cells.Item($row, 11).Value = 99.42
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}

Write-Output "Applied trade #27 close update"
